$wb = $excel.ActiveWorkbook

# The current active sheet is "MegaMenuInfo" (2nd sheet). Its selection will
# move from A5 to B1:C1 and it will no longer be the selected tab once the
# new sheet is appended and activated.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B1:C1").Select()

# Add the new "TimingScroll" worksheet as the last tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newws = $wb.Worksheets.Add($null, $lastSheet)
$newws.Name = "TimingScroll"

# Header row - reuse the existing "Path"/"ContentType" header formatting
# from the MegaMenuInfo sheet (bold, shaded style).
$newws.Range("A1").Value = "Path"
$newws.Range("B1").Value = "ContentType"
$ws2.Range("A1:B1").Copy()
$newws.Range("A1:B1").PasteSpecial(-4122)

# Data rows.
$newws.Range("A2").Value = "/news-events"
$newws.Range("B2").Value = "Landing"

$newws.Range("A3").Value = "/espanol/tipos/vesicula-biliar"
$newws.Range("B3").Value = "CTHP"

$newws.Range("A4").Value = "/about-cancer/treatment/types/immunotherapy"
$newws.Range("B4").Value = "Article"

# Best-fit the columns to their content.
$newws.Columns.Item(1).EntireColumn.AutoFit()
$newws.Columns.Item(2).EntireColumn.AutoFit()

# Leave the selection on A5 on the new, now-active sheet.
$newws.Range("A5").Select()
